$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = 111866919
$ws.Range("B11").Value = 90689
$ws.Range("E11").Value = 5966
$ws.Range("F11").Value = "Motaggsvamp"
$ws.Range("G11").Value = "Sarcodon squamosus"
$ws.Range("H11").Value = "(Schaeff.) Quél."
$ws.Range("Q11").Value = 703070.5942336121
$ws.Range("R11").Value = 7299535.948440861

# Row 12
$ws.Range("A12").Value = 111867661
$ws.Range("B12").Value = 73692
$ws.Range("E12").Value = 310
$ws.Range("F12").Value = "Nordlig nållav"
$ws.Range("G12").Value = "Chaenotheca laevigata"
$ws.Range("H12").Value = "Nádv."
$ws.Range("Q12").Value = 703308.4646664646
$ws.Range("R12").Value = 7299302.011735545

# Row 13
$ws.Range("A13").Value = 111867696
$ws.Range("B13").Value = 90660
$ws.Range("E13").Value = 4362
$ws.Range("F13").Value = "Blå taggsvamp"
$ws.Range("G13").Value = "Hydnellum caeruleum"
$ws.Range("H13").Value = "(Hornem.) P.Karst."
$ws.Range("Q13").Value = 703310.8095286442
$ws.Range("R13").Value = 7299298.053094583

# Row 17
$ws.Range("A17").Value = 111867007
$ws.Range("B17").Value = 90709
$ws.Range("E17").Value = 5448
$ws.Range("F17").Value = "Svartvit taggsvamp"
$ws.Range("G17").Value = "Phellodon connatus"
$ws.Range("H17").Value = "(Schultz) nom.prov"
$ws.Range("Q17").Value = 703070.0396593859
$ws.Range("R17").Value = 7299502.915184345

# Row 18
$ws.Range("A18").Value = 111867419
$ws.Range("B18").Value = 90658
$ws.Range("E18").Value = 4361
$ws.Range("F18").Value = "Orange taggsvamp"
$ws.Range("G18").Value = "Hydnellum aurantiacum"
$ws.Range("H18").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q18").Value = 703159.5881134692
$ws.Range("R18").Value = 7299375.205745419

# Row 36
$ws.Range("A36").Value = 111867075
$ws.Range("B36").Value = 90660
$ws.Range("E36").Value = 4362
$ws.Range("F36").Value = "Blå taggsvamp"
$ws.Range("G36").Value = "Hydnellum caeruleum"
$ws.Range("H36").Value = "(Hornem.) P.Karst."
$ws.Range("Q36").Value = 703079.296544011
$ws.Range("R36").Value = 7299482.94824858

# Row 37
$ws.Range("A37").Value = 111866994
$ws.Range("B37").Value = 90709
$ws.Range("E37").Value = 5448
$ws.Range("F37").Value = "Svartvit taggsvamp"
$ws.Range("G37").Value = "Phellodon connatus"
$ws.Range("H37").Value = "(Schultz) nom.prov"
$ws.Range("Q37").Value = 703114.8550411762
$ws.Range("R37").Value = 7299511.445840456

# Row 38
$ws.Range("A38").Value = 111867403
$ws.Range("B38").Value = 90660
$ws.Range("E38").Value = 4362
$ws.Range("F38").Value = "Blå taggsvamp"
$ws.Range("G38").Value = "Hydnellum caeruleum"
$ws.Range("H38").Value = "(Hornem.) P.Karst."
$ws.Range("Q38").Value = 703140.5813816102
$ws.Range("R38").Value = 7299387.059685718
